$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: clear the error-message cell (C2)
$ws.Range("C2").ClearContents()

# Row 3: B3 numeric 123456 becomes a text value "123456" (quote-prefixed),
# and its error-message cell (C3) is cleared
$ws.Range("B3").Value = "'123456"
$ws.Range("C3").ClearContents()

# Row 4: B4 numeric 123456 becomes a text value "123456" (quote-prefixed)
$ws.Range("B4").Value = "'123456"

# Row 5: B5 numeric 1234567 becomes a text value "1234567" (quote-prefixed)
$ws.Range("B5").Value = "'1234567"

# Row 4 & 5: error-message cells are replaced with the new unified message
$ws.Range("C4").Value = "Ensure valid username/password!"
$ws.Range("C5").Value = "Ensure valid username/password!"

# Update the selected/active cell shown in the saved view
$ws.Range("C11").Select()
